$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

$ws.Range("B1").Value = "Minnesota"

$c1 = $ws.Range("C1")
$c1.Value = 44509
$c1.NumberFormat = "mm-dd-yy"
